# Update row 2 of test_data/make_payment.xlsx:
#  - A2 username token
#  - B2 phone number
#  - C2 display name
#  - E2 payment date (kept as literal text, matching the source column header
#    "Payment Date (YYYY-MM-DD)")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 (leading "+") and E2 (date-looking string) must stay plain text so
# Excel doesn't reinterpret them as a number / date serial.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "677f5c84fc34afaac4ae"
$ws.Range("B2").Value = "+74267426016"
$ws.Range("C2").Value = "Automation User 16"
$ws.Range("E2").Value = "2026-01-05"
